$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 142, shifting rows 142:234 down to 143:235
$ws.Rows.Item(142).Insert()

# Fill the newly inserted row 142 with the new data record.
# (Same constant fields as the surrounding rows, new measurement values.)
$ws.Cells.Item(142, 1).Value = 4
$ws.Cells.Item(142, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(142, 3).Value = "Los Lagos"
$ws.Cells.Item(142, 4).Value = 44603
$ws.Cells.Item(142, 4).Style = $ws.Cells.Item(143, 4).Style
$ws.Cells.Item(142, 4).NumberFormat = $ws.Cells.Item(143, 4).NumberFormat
$ws.Cells.Item(142, 5).Value = 10
$ws.Cells.Item(142, 6).Value = 100112037
$ws.Cells.Item(142, 7).Value = "Cebollín"
$ws.Cells.Item(142, 8).Value = "Sin especificar"
$ws.Cells.Item(142, 9).Value = "Primera"
$ws.Cells.Item(142, 10).Value = 160
$ws.Cells.Item(142, 11).Value = 5500
$ws.Cells.Item(142, 12).Value = 6000
$ws.Cells.Item(142, 13).Value = 5750
$ws.Cells.Item(142, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(142, 15).Value = "Región Metropolitana"
$ws.Cells.Item(142, 16).Value = 160
$ws.Cells.Item(142, 17).Value = 36
$ws.Cells.Item(142, 18).Value = "Hortaliza"
